# Applies updated market-board pricing figures to the per-job Leve profit
# sheets (columns H-N) as captured by the scheduled pricing-update run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 92: Enchanted Koppranickel Ink
$ws.Range("H92").Value = 38462716
$ws.Range("I92").Value = 45455700
$ws.Range("K92").Value = 45455700
$ws.Range("M92").Value = -45454452

# Row 100: Beetle Glue
$ws.Range("H100").Value = 4666.6665
$ws.Range("I100").Value = 4000
$ws.Range("J100").Value = 4750
$ws.Range("K100").Value = 4000
$ws.Range("L100").Value = 4750
$ws.Range("M100").Value = -3459
$ws.Range("N100").Value = -5832

# Row 112: Superior Spiritbond Potion
$ws.Range("H112").Value = 1152089.1
$ws.Range("I112").Value = 2633.125
$ws.Range("J112").Value = 1589977.2
$ws.Range("K112").Value = 7899.375
$ws.Range("L112").Value = 4769931.6
$ws.Range("M112").Value = -6791.375
$ws.Range("N112").Value = -4772147.6

# Row 116: Growth Formula Kappa
$ws.Range("H116").Value = 17255.928
$ws.Range("I116").Value = 38278.8
$ws.Range("J116").Value = 5576.5557
$ws.Range("K116").Value = 38278.8
$ws.Range("L116").Value = 5576.5557
$ws.Range("M116").Value = -34836.8
$ws.Range("N116").Value = -12460.5557

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Steel Ingot
$ws.Range("H32").Value = 24419468
$ws.Range("I32").Value = 27172794
$ws.Range("J32").Value = 8932011
$ws.Range("K32").Value = 27172794
$ws.Range("L32").Value = 8932011
$ws.Range("M32").Value = -27172507
$ws.Range("N32").Value = -8932585

# Row 45: Mythril Ingot
$ws.Range("H45").Value = 4585.727
$ws.Range("I45").Value = 4828.1113
$ws.Range("K45").Value = 4828.1113
$ws.Range("M45").Value = -4451.1113

# Row 112: Deepgold Gloves of Fending
$ws.Range("H112").Value = 51989
$ws.Range("J112").Value = 51989
$ws.Range("L112").Value = 51989
$ws.Range("N112").Value = -54943

# Row 114: Bluespirit Gauntlets of Fending
$ws.Range("H114").Value = 59138
$ws.Range("J114").Value = 59138
$ws.Range("L114").Value = 59138
$ws.Range("N114").Value = -67816

# Row 119: Dwarven Mythril Chainmail of Fending
$ws.Range("H119").Value = 80698
$ws.Range("J119").Value = 80698
$ws.Range("L119").Value = 80698
$ws.Range("N119").Value = -90374

# Row 122: High Durium Nugget
$ws.Range("H122").Value = 5194.8423
$ws.Range("I122").Value = 3870
$ws.Range("K122").Value = 11610
$ws.Range("M122").Value = -9160

# Row 125: High Durium Armor of Fending
$ws.Range("H125").Value = 51467.832
$ws.Range("J125").Value = 51467.832
$ws.Range("L125").Value = 51467.832
$ws.Range("N125").Value = -61307.832

$ws = $wb.Worksheets.Item("BSM")
# Row 134: Ruthenium Ingot
$ws.Range("H134").Value = 2383194.2
$ws.Range("I134").Value = 2647386.2
$ws.Range("K134").Value = 7942158.600000001
$ws.Range("M134").Value = -7939623.600000001

$ws = $wb.Worksheets.Item("CRP")
# Row 2: Bone Harpoon
$ws.Range("H2").Value = 1673.625
$ws.Range("I2").Value = 1627
$ws.Range("J2").Value = 2000
$ws.Range("K2").Value = 1627
$ws.Range("L2").Value = 2000
$ws.Range("M2").Value = -1514
$ws.Range("N2").Value = -2226

# Row 31: Walnut Lumber
$ws.Range("H31").Value = 7762.409
$ws.Range("J31").Value = 8238.75
$ws.Range("L31").Value = 8238.75
$ws.Range("N31").Value = -8828.75

# Row 34: Walnut Lumber
$ws.Range("H34").Value = 7762.409
$ws.Range("J34").Value = 8238.75
$ws.Range("L34").Value = 8238.75
$ws.Range("N34").Value = -8642.75

# Row 108: White Oak Fishing Rod
$ws.Range("H108").Value = 164000
$ws.Range("J108").Value = 164000
$ws.Range("L108").Value = 164000
$ws.Range("N108").Value = -171680

$ws = $wb.Worksheets.Item("CUL")
# Row 12: Kukuru Butter
$ws.Range("H12").Value = 1111230.8
$ws.Range("I12").Value = 69
$ws.Range("J12").Value = 1428705.6
$ws.Range("K12").Value = 207
$ws.Range("L12").Value = 4286116.800000001
$ws.Range("M12").Value = -34
$ws.Range("N12").Value = -4286462.800000001

# Row 19: Parsnip Salad
$ws.Range("H19").Value = 1000
$ws.Range("J19").Value = 1000
$ws.Range("L19").Value = 3000
$ws.Range("N19").Value = -3348

# Row 107: Frantoio Oil
$ws.Range("H107").Value = 809.2059
$ws.Range("I107").Value = 1162.8572
$ws.Range("J107").Value = 561.65
$ws.Range("K107").Value = 3488.5716
$ws.Range("L107").Value = 1684.95
$ws.Range("M107").Value = -1568.5716
$ws.Range("N107").Value = -5524.95

# Row 117: Peppered Popotoes
$ws.Range("H117").Value = 3094.125
$ws.Range("J117").Value = 4184.8
$ws.Range("L117").Value = 12554.4
$ws.Range("N117").Value = -19438.4

# Row 120: Paella
$ws.Range("H120").Value = 22503
$ws.Range("I120").Value = 30
$ws.Range("K120").Value = 90
$ws.Range("M120").Value = 4748

# Row 132: Cooking Mezcal
$ws.Range("H132").Value = 1043.6471
$ws.Range("J132").Value = 1485.3334
$ws.Range("L132").Value = 13368.0006
$ws.Range("N132").Value = -18428.0006

$ws = $wb.Worksheets.Item("GSM")
# Row 41: Worm Fang Needle
$ws.Range("H41").Value = 1048.6666
$ws.Range("I41").Value = 1048.6666
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 1048.6666
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -693.6666
$ws.Range("N41").ClearContents()

# Row 80: Hardsilver Ingot
$ws.Range("H80").Value = 3133.2222
$ws.Range("I80").Value = 2928.4285
$ws.Range("K80").Value = 2928.4285
$ws.Range("M80").Value = -1930.4285

# Row 83: Hardsilver Ingot
$ws.Range("H83").Value = 3133.2222
$ws.Range("I83").Value = 2928.4285
$ws.Range("K83").Value = 14642.1425
$ws.Range("M83").Value = -9650.1425

# Row 123: Ametrine Ring of Fending
$ws.Range("H123").Value = 29888
$ws.Range("J123").Value = 29888
$ws.Range("L123").Value = 29888
$ws.Range("N123").Value = -34788

# Row 124: Pewter Pendulums
$ws.Range("H124").Value = 114943.75
$ws.Range("J124").Value = 114943.75
$ws.Range("L124").Value = 114943.75
$ws.Range("N124").Value = -124763.75

# Row 132: Lar Ingot
$ws.Range("H132").Value = 3865.8462
$ws.Range("I132").Value = 3865.8462
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 11597.5386
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -9067.5386
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 22: Aldgoat Leather
$ws.Range("H22").Value = 3151.2144
$ws.Range("I22").Value = 1934.7778
$ws.Range("J22").Value = 5340.8
$ws.Range("K22").Value = 1934.7778
$ws.Range("L22").Value = 5340.8
$ws.Range("M22").Value = -1639.7778
$ws.Range("N22").Value = -5930.8

# Row 27: Aldgoat Leather
$ws.Range("H27").Value = 3151.2144
$ws.Range("I27").Value = 1934.7778
$ws.Range("J27").Value = 5340.8
$ws.Range("K27").Value = 1934.7778
$ws.Range("L27").Value = 5340.8
$ws.Range("M27").Value = -1827.7778
$ws.Range("N27").Value = -5554.8

# Row 32: Goatskin Targe
$ws.Range("H32").Value = 24499.666
$ws.Range("I32").Value = 23399.8
$ws.Range("K32").Value = 23399.8
$ws.Range("M32").Value = -23082.8

# Row 108: Smilodonskin Trousers of Maiming
$ws.Range("H108").Value = 35282.332
$ws.Range("J108").Value = 35282.332
$ws.Range("L108").Value = 35282.332
$ws.Range("N108").Value = -42962.332

# Row 132: Silver Lobo Leather
$ws.Range("H132").Value = 4178.6206
$ws.Range("I132").Value = 4026.2173
$ws.Range("J132").Value = 4762.8335
$ws.Range("K132").Value = 12078.6519
$ws.Range("L132").Value = 14288.5005
$ws.Range("M132").Value = -9548.651899999999
$ws.Range("N132").Value = -19348.5005

$ws = $wb.Worksheets.Item("WVR")
# Row 122: Dark Hempen Cloth
$ws.Range("H122").Value = 13555
$ws.Range("I122").Value = 13577.6
$ws.Range("K122").Value = 40732.8
$ws.Range("M122").Value = -38282.8

# Row 136: Sarcenet Cloth
$ws.Range("H136").Value = 11766510
$ws.Range("I136").Value = 1758.2307
$ws.Range("J136").Value = 50001950
$ws.Range("K136").Value = 5274.6921
$ws.Range("L136").Value = 150005850
$ws.Range("M136").Value = -2724.6921
$ws.Range("N136").Value = -150010950

# Row 140: Thunderyards Silk Gloves of Casting
$ws.Range("H140").Value = 56107
$ws.Range("J140").Value = 56107
$ws.Range("L140").Value = 56107
$ws.Range("N140").Value = -66467
